$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1, columns A-L currently hold "Comp.1".."Comp.12";
# rename them to "PC1".."PC12".
for ($i = 1; $i -le 12; $i++) {
    $cell = $ws.Cells.Item(1, $i)
    $cell.Value = "PC$i"
}
